# Add a new worksheet "apiTest" after the existing "LoginTest" sheet,
# populate it with the API registration test data, add the mailto
# hyperlink for the email cell, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "apiTest"

# Header row
$newSheet.Range("A1").Value = "email"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("C1").Value = "firstName"
$newSheet.Range("D1").Value = "lastName"

# Data row
$newSheet.Range("A2").Value = "srdjan.rados@htecgroup.com"
$newSheet.Range("B2").Value = "Qwertysha1@"
$newSheet.Range("C2").Value = "Srdjan"
$newSheet.Range("D2").Value = "Rados"

# Hyperlink on the email cell
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:srdjan.rados@htecgroup.com", "", "", "srdjan.rados@htecgroup.com")

# Make the new sheet the active / selected sheet
$newSheet.Activate()
